$d = $word.ActiveDocument

$d.Content.Find.Execute("12-8=4", $true, $false, $false, $false, $false, $true, 1, $false, "78-67=11", 2) | Out-Null
$d.Content.Find.Execute("23+53=76", $true, $false, $false, $false, $false, $true, 1, $false, "22+32=54", 2) | Out-Null
$d.Content.Find.Execute("84-17=67", $true, $false, $false, $false, $false, $true, 1, $false, "13+54=67", 2) | Out-Null
$d.Content.Find.Execute("7+82=89", $true, $false, $false, $false, $false, $true, 1, $false, "81-11=70", 2) | Out-Null
$d.Content.Find.Execute("95-0=95", $true, $false, $false, $false, $false, $true, 1, $false, "79+3=82", 2) | Out-Null
$d.Content.Find.Execute("95-31=64", $true, $false, $false, $false, $false, $true, 1, $false, "96-31=65", 2) | Out-Null
$d.Content.Find.Execute("28+65=93", $true, $false, $false, $false, $false, $true, 1, $false, "46+20=66", 2) | Out-Null
$d.Content.Find.Execute("93+0=93", $true, $false, $false, $false, $false, $true, 1, $false, "93-56=37", 2) | Out-Null
$d.Content.Find.Execute("79-27=52", $true, $false, $false, $false, $false, $true, 1, $false, "51-14=37", 2) | Out-Null
$d.Content.Find.Execute("70-23=47", $true, $false, $false, $false, $false, $true, 1, $false, "26+0=26", 2) | Out-Null
$d.Content.Find.Execute("98-0=98", $true, $false, $false, $false, $false, $true, 1, $false, "2+51=53", 2) | Out-Null
$d.Content.Find.Execute("12+77=89", $true, $false, $false, $false, $false, $true, 1, $false, "69-66=3", 2) | Out-Null
$d.Content.Find.Execute("48+26=74", $true, $false, $false, $false, $false, $true, 1, $false, "25+1=26", 2) | Out-Null
$d.Content.Find.Execute("27+1=28", $true, $false, $false, $false, $false, $true, 1, $false, "85-54=31", 2) | Out-Null
$d.Content.Find.Execute("15+13=28", $true, $false, $false, $false, $false, $true, 1, $false, "29+32=61", 2) | Out-Null
$d.Content.Find.Execute("20+31=51", $true, $false, $false, $false, $false, $true, 1, $false, "5+66=71", 2) | Out-Null
$d.Content.Find.Execute("25+74=99", $true, $false, $false, $false, $false, $true, 1, $false, "29+16=45", 2) | Out-Null
$d.Content.Find.Execute("49+38=87", $true, $false, $false, $false, $false, $true, 1, $false, "58-54=4", 2) | Out-Null
$d.Content.Find.Execute("31-25=6", $true, $false, $false, $false, $false, $true, 1, $false, "31+5=36", 2) | Out-Null
$d.Content.Find.Execute("96-69=27", $true, $false, $false, $false, $false, $true, 1, $false, "70+14=84", 2) | Out-Null
$d.Content.Find.Execute("83+11=94", $true, $false, $false, $false, $false, $true, 1, $false, "8+28=36", 2) | Out-Null
$d.Content.Find.Execute("8+67=75", $true, $false, $false, $false, $false, $true, 1, $false, "15+2=17", 2) | Out-Null
$d.Content.Find.Execute("37-32=5", $true, $false, $false, $false, $false, $true, 1, $false, "30-27=3", 2) | Out-Null
$d.Content.Find.Execute("55-31=24", $true, $false, $false, $false, $false, $true, 1, $false, "49+0=49", 2) | Out-Null
$d.Content.Find.Execute("85-70=15", $true, $false, $false, $false, $false, $true, 1, $false, "93+6=99", 2) | Out-Null
$d.Content.Find.Execute("18-11=7", $true, $false, $false, $false, $false, $true, 1, $false, "16-6=10", 2) | Out-Null
$d.Content.Find.Execute("0+64=64", $true, $false, $false, $false, $false, $true, 1, $false, "91-69=22", 2) | Out-Null
$d.Content.Find.Execute("39+37=76", $true, $false, $false, $false, $false, $true, 1, $false, "92-80=12", 2) | Out-Null
$d.Content.Find.Execute("81-49=32", $true, $false, $false, $false, $false, $true, 1, $false, "12+13=25", 2) | Out-Null
$d.Content.Find.Execute("53-26=27", $true, $false, $false, $false, $false, $true, 1, $false, "48+30=78", 2) | Out-Null
$d.Content.Find.Execute("27-9=18", $true, $false, $false, $false, $false, $true, 1, $false, "39-25=14", 2) | Out-Null
$d.Content.Find.Execute("52+35=87", $true, $false, $false, $false, $false, $true, 1, $false, "18+55=73", 2) | Out-Null
$d.Content.Find.Execute("22+20=42", $true, $false, $false, $false, $false, $true, 1, $false, "84-42=42", 2) | Out-Null
$d.Content.Find.Execute("48+21=69", $true, $false, $false, $false, $false, $true, 1, $false, "70-10=60", 2) | Out-Null
$d.Content.Find.Execute("65-12=53", $true, $false, $false, $false, $false, $true, 1, $false, "55-41=14", 2) | Out-Null
$d.Content.Find.Execute("8+8=16", $true, $false, $false, $false, $false, $true, 1, $false, "5+83=88", 2) | Out-Null
$d.Content.Find.Execute("13+63=76", $true, $false, $false, $false, $false, $true, 1, $false, "78+4=82", 2) | Out-Null
$d.Content.Find.Execute("6-2=4", $true, $false, $false, $false, $false, $true, 1, $false, "15+36=51", 2) | Out-Null
$d.Content.Find.Execute("77+12=89", $true, $false, $false, $false, $false, $true, 1, $false, "89-24=65", 2) | Out-Null
$d.Content.Find.Execute("23+67=90", $true, $false, $false, $false, $false, $true, 1, $false, "88-7=81", 2) | Out-Null
$d.Content.Find.Execute("54-5=49", $true, $false, $false, $false, $false, $true, 1, $false, "6+82=88", 2) | Out-Null
$d.Content.Find.Execute("82-58=24", $true, $false, $false, $false, $false, $true, 1, $false, "13+46=59", 2) | Out-Null
$d.Content.Find.Execute("53-29=24", $true, $false, $false, $false, $false, $true, 1, $false, "89-69=20", 2) | Out-Null
$d.Content.Find.Execute("52-1=51", $true, $false, $false, $false, $false, $true, 1, $false, "38+50=88", 2) | Out-Null
$d.Content.Find.Execute("38+46=84", $true, $false, $false, $false, $false, $true, 1, $false, "76-26=50", 2) | Out-Null
$d.Content.Find.Execute("76-37=39", $true, $false, $false, $false, $false, $true, 1, $false, "38+54=92", 2) | Out-Null
$d.Content.Find.Execute("38-28=10", $true, $false, $false, $false, $false, $true, 1, $false, "92+1=93", 2) | Out-Null
$d.Content.Find.Execute("81+2=83", $true, $false, $false, $false, $false, $true, 1, $false, "12-11=1", 2) | Out-Null
$d.Content.Find.Execute("51+34=85", $true, $false, $false, $false, $false, $true, 1, $false, "38-6=32", 2) | Out-Null
$d.Content.Find.Execute("35+48=83", $true, $false, $false, $false, $false, $true, 1, $false, "42-35=7", 2) | Out-Null
$d.Content.Find.Execute("85-31=54", $true, $false, $false, $false, $false, $true, 1, $false, "33-25=8", 2) | Out-Null
$d.Content.Find.Execute("70+4=74", $true, $false, $false, $false, $false, $true, 1, $false, "81-22=59", 2) | Out-Null
$d.Content.Find.Execute("76-64=12", $true, $false, $false, $false, $false, $true, 1, $false, "13+8=21", 2) | Out-Null
$d.Content.Find.Execute("91-25=66", $true, $false, $false, $false, $false, $true, 1, $false, "17+20=37", 2) | Out-Null
$d.Content.Find.Execute("20+27=47", $true, $false, $false, $false, $false, $true, 1, $false, "90-35=55", 2) | Out-Null
$d.Content.Find.Execute("73-56=17", $true, $false, $false, $false, $false, $true, 1, $false, "29-1=28", 2) | Out-Null
$d.Content.Find.Execute("56-28=28", $true, $false, $false, $false, $false, $true, 1, $false, "60-53=7", 2) | Out-Null
$d.Content.Find.Execute("91-32=59", $true, $false, $false, $false, $false, $true, 1, $false, "30-20=10", 2) | Out-Null
$d.Content.Find.Execute("28-4=24", $true, $false, $false, $false, $false, $true, 1, $false, "54-39=15", 2) | Out-Null
$d.Content.Find.Execute("43-28=15", $true, $false, $false, $false, $false, $true, 1, $false, "76-76=0", 2) | Out-Null
$d.Content.Find.Execute("63+8=71", $true, $false, $false, $false, $false, $true, 1, $false, "52+16=68", 2) | Out-Null
$d.Content.Find.Execute("66-43=23", $true, $false, $false, $false, $false, $true, 1, $false, "74-35=39", 2) | Out-Null
$d.Content.Find.Execute("11+54=65", $true, $false, $false, $false, $false, $true, 1, $false, "33+24=57", 2) | Out-Null
$d.Content.Find.Execute("38-36=2", $true, $false, $false, $false, $false, $true, 1, $false, "20+26=46", 2) | Out-Null
$d.Content.Find.Execute("9+2=11", $true, $false, $false, $false, $false, $true, 1, $false, "72-28=44", 2) | Out-Null
$d.Content.Find.Execute("54-1=53", $true, $false, $false, $false, $false, $true, 1, $false, "57+2=59", 2) | Out-Null
$d.Content.Find.Execute("52-40=12", $true, $false, $false, $false, $false, $true, 1, $false, "77-52=25", 2) | Out-Null
$d.Content.Find.Execute("52-14=38", $true, $false, $false, $false, $false, $true, 1, $false, "2+13=15", 2) | Out-Null
$d.Content.Find.Execute("64+28=92", $true, $false, $false, $false, $false, $true, 1, $false, "43+42=85", 2) | Out-Null
$d.Content.Find.Execute("73-10=63", $true, $false, $false, $false, $false, $true, 1, $false, "59-34=25", 2) | Out-Null
$d.Content.Find.Execute("66-38=28", $true, $false, $false, $false, $false, $true, 1, $false, "59+27=86", 2) | Out-Null
$d.Content.Find.Execute("11+72=83", $true, $false, $false, $false, $false, $true, 1, $false, "7+76=83", 2) | Out-Null
$d.Content.Find.Execute("49-7=42", $true, $false, $false, $false, $false, $true, 1, $false, "92-43=49", 2) | Out-Null
$d.Content.Find.Execute("89-25=64", $true, $false, $false, $false, $false, $true, 1, $false, "38+33=71", 2) | Out-Null
$d.Content.Find.Execute("9-4=5", $true, $false, $false, $false, $false, $true, 1, $false, "92-13=79", 2) | Out-Null
$d.Content.Find.Execute("56+2=58", $true, $false, $false, $false, $false, $true, 1, $false, "98-78=20", 2) | Out-Null
$d.Content.Find.Execute("97-2=95", $true, $false, $false, $false, $false, $true, 1, $false, "34-9=25", 2) | Out-Null
$d.Content.Find.Execute("51+0=51", $true, $false, $false, $false, $false, $true, 1, $false, "91-63=28", 2) | Out-Null
$d.Content.Find.Execute("3+66=69", $true, $false, $false, $false, $false, $true, 1, $false, "77+18=95", 2) | Out-Null
$d.Content.Find.Execute("73-69=4", $true, $false, $false, $false, $false, $true, 1, $false, "89+5=94", 2) | Out-Null
$d.Content.Find.Execute("15+60=75", $true, $false, $false, $false, $false, $true, 1, $false, "23+6=29", 2) | Out-Null
$d.Content.Find.Execute("92-72=20", $true, $false, $false, $false, $false, $true, 1, $false, "44+41=85", 2) | Out-Null
$d.Content.Find.Execute("1+18=19", $true, $false, $false, $false, $false, $true, 1, $false, "75-65=10", 2) | Out-Null
$d.Content.Find.Execute("39-38=1", $true, $false, $false, $false, $false, $true, 1, $false, "95-68=27", 2) | Out-Null
$d.Content.Find.Execute("24+8=32", $true, $false, $false, $false, $false, $true, 1, $false, "64+12=76", 2) | Out-Null
$d.Content.Find.Execute("74-29=45", $true, $false, $false, $false, $false, $true, 1, $false, "43+17=60", 2) | Out-Null
$d.Content.Find.Execute("2-0=2", $true, $false, $false, $false, $false, $true, 1, $false, "56+3=59", 2) | Out-Null
$d.Content.Find.Execute("68-30=38", $true, $false, $false, $false, $false, $true, 1, $false, "62-14=48", 2) | Out-Null
$d.Content.Find.Execute("26+70=96", $true, $false, $false, $false, $false, $true, 1, $false, "65-11=54", 2) | Out-Null
$d.Content.Find.Execute("5+39=44", $true, $false, $false, $false, $false, $true, 1, $false, "24+44=68", 2) | Out-Null
$d.Content.Find.Execute("31+45=76", $true, $false, $false, $false, $false, $true, 1, $false, "17+74=91", 2) | Out-Null
$d.Content.Find.Execute("57-36=21", $true, $false, $false, $false, $false, $true, 1, $false, "67+6=73", 2) | Out-Null
$d.Content.Find.Execute("71-44=27", $true, $false, $false, $false, $false, $true, 1, $false, "62-1=61", 2) | Out-Null
$d.Content.Find.Execute("42-38=4", $true, $false, $false, $false, $false, $true, 1, $false, "65-24=41", 2) | Out-Null
$d.Content.Find.Execute("94-71=23", $true, $false, $false, $false, $false, $true, 1, $false, "86-64=22", 2) | Out-Null
$d.Content.Find.Execute("91-61=30", $true, $false, $false, $false, $false, $true, 1, $false, "34-22=12", 2) | Out-Null
$d.Content.Find.Execute("34+37=71", $true, $false, $false, $false, $false, $true, 1, $false, "9+61=70", 2) | Out-Null
$d.Content.Find.Execute("60+28=88", $true, $false, $false, $false, $false, $true, 1, $false, "97-22=75", 2) | Out-Null
$d.Content.Find.Execute("48+17=65", $true, $false, $false, $false, $false, $true, 1, $false, "33+14=47", 2) | Out-Null
$d.Content.Find.Execute("81-5=76", $true, $false, $false, $false, $false, $true, 1, $false, "5+8=13", 2) | Out-Null

Write-Output "Replaced 100 entries"
